# Forward Look workbook update
# - Update the "as at" date in the intro paragraph (12 June 2023 -> 22 June 2023)
#   and wrap it onto a second, indented line.
# - Insert a new "Type" column (F) populated with "standard" for every
#   publication row, mirroring the formatting already used in column D/E.
# - Resize a few columns (A, C, D, E, F).
# - Extend the banded-row conditional formatting to cover the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the descriptive paragraph in A2
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics`n                                            that have been pre-announced on the gov.uk release calendar as at 22 June 2023"
# The embedded line break otherwise makes the engine apply a custom row
# height; re-fit the row so it keeps using the sheet's default height, just
# like the original (unmodified) row did.
$ws.Rows.Item(2).AutoFit()

# ---------------------------------------------------------------------------
# 2. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.83   # A: 18.71 -> 14.71
$ws.Columns.Item(3).ColumnWidth = 29.83   # C: 24.71 -> 30.71
$ws.Columns.Item(4).ColumnWidth = 9.83    # D: 12.71 -> 10.71
$ws.Columns.Item(5).ColumnWidth = 9.83    # E: 12.71 -> 10.71 (stays hidden)
$ws.Columns.Item(6).ColumnWidth = 13.83   # F: new column, width 14.71

# ---------------------------------------------------------------------------
# 3. New column F ("Type") - copy formatting from column D so that styled
#    (s="5") and unstyled rows line up exactly like the rest of the table,
#    then fill in the header and the "standard" values.
# ---------------------------------------------------------------------------
$ws.Range("D4").Copy()
$ws.Range("F4").PasteSpecial(-4122)

$ws.Range("D5:D61").Copy()
$ws.Range("F5:F61").PasteSpecial(-4122)

$ws.Range("F4").Value = "Type"

$ws.Range("F5:F7").Value = "standard"
$ws.Range("F9:F19").Value = "standard"
$ws.Range("F21:F24").Value = "standard"
$ws.Range("F27:F28").Value = "standard"
$ws.Range("F30:F33").Value = "standard"
$ws.Range("F35:F42").Value = "standard"
$ws.Range("F44:F53").Value = "standard"
$ws.Range("F58:F61").Value = "standard"

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Extend the banded-row conditional formatting (previously A5:E61) so it
#    also covers the new column F.
# ---------------------------------------------------------------------------
$cfs = $ws.Range("A5:E61").FormatConditions
for ($i = 1; $i -le $cfs.Count; $i++) {
    $fc = $cfs.Item($i)
    if ($fc.AppliesTo.Address() -eq '$A$5:$E$61') {
        $fc.ModifyAppliesToRange($ws.Range("A5:F61"))
    }
}
